$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.274.32"
$ws.Range("E2").Value = "  +2.01%  "

$ws.Range("D3").Value = "1.815.34"
$ws.Range("E3").Value = "  +3.49%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.87"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4380"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3670"
$ws.Range("E8").Value = "  +0.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.81"
$ws.Range("E9").Value = "  -0.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07673"
$ws.Range("E10").Value = "  +2.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.140"
$ws.Range("E11").Value = "  +1.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9998"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.99"
$ws.Range("E13").Value = "  +1.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.309"
$ws.Range("E14").Value = "  +2.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.489"
$ws.Range("E15").Value = "  +3.24%  "

$ws.Range("D16").Value = "1.819.04"
$ws.Range("E16").Value = "  +3.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.03"
$ws.Range("E17").Value = "  +7.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001079"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06470"
$ws.Range("E19").Value = "  +4.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9993"
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.38"
$ws.Range("E21").Value = "  +1.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.240"
$ws.Range("E22").Value = "  +1.25%  "

$ws.Range("D23").Value = "28.292.62"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.56"
$ws.Range("E24").Value = "  -0.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.127"
$ws.Range("E25").Value = "  -8.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.03"
$ws.Range("E26").Value = "  +5.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.72"
$ws.Range("E27").Value = "  +0.60%  "

$ws.Range("D28").Value = "2.026.75"
$ws.Range("E28").Value = "  +3.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.276"
$ws.Range("E29").Value = "  -3.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.49"
$ws.Range("E30").Value = "  +1.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.212"
$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.991"
$ws.Range("E32").Value = "  +4.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09142"
$ws.Range("E33").Value = "  -0.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.539"
$ws.Range("E34").Value = "  -3.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.03"
$ws.Range("E35").Value = "  +2.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02362"
$ws.Range("E36").Value = "  +2.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.243"
$ws.Range("E37").Value = "  +2.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2172"
$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6600"
$ws.Range("E39").Value = "  +1.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06208"
$ws.Range("E40").Value = "  +1.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.199"
$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.058"
$ws.Range("E42").Value = "  +1.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.430"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9989"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.78"
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6102"
$ws.Range("E46").Value = "  +2.53%  "

$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.021"
$ws.Range("E48").Value = "  +2.06%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.32"
$ws.Range("E49").Value = "  -0.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.161"
$ws.Range("E50").Value = "  +3.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06994"
$ws.Range("E51").Value = "  +1.37%  "

